$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.488.21'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.332.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '546.94'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.13'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.331.49'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.51'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.150'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.335'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.51'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.40%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.747.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.394.99'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.329.73'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.57'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.98'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.07'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.58'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.93'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.42%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.88'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.39'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.22'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.79'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.73'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0732'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.92'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.37'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.382'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.93'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("B37").Value = 'USDe'
$ws.Range("C37").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.11'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '319.56'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.08'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.53'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.72'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.48'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0942'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.11'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.571'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0496'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0216'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0220'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.98'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.48%  '
